$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.045.17"
$ws.Range("E2").Value = "  -0.92%  "

# Row 3
$ws.Range("D3").Value = "2.198.50"
$ws.Range("E3").Value = "  -2.32%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.97%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("E9").Value = "  -8.32%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.77%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.16%  "

# Row 13
$ws.Range("E13").Value = "  -1.04%  "

# Row 14
$ws.Range("D14").Value = "2.533.15"
$ws.Range("E14").Value = "  -2.16%  "

# Row 15
$ws.Range("D15").Value = "2.260.01"
$ws.Range("E15").Value = "  -5.50%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.00%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.779"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.68%  "

# Row 18
$ws.Range("D18").Value = "43.734.34"
$ws.Range("E18").Value = "  -0.91%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0900"
$ws.Range("E19").Value = "  -7.21%  "

# Row 20
$ws.Range("E20").Value = "  -9.33%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -12.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "63.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -14.19%  "

# Row 25
$ws.Range("E25").Value = "  +0.91%  "

# Row 26
$ws.Range("E26").Value = "  -8.66%  "

# Row 27
$ws.Range("E27").Value = "  +0.08%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.55%  "

# Row 29
$ws.Range("E29").Value = "  -6.42%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "148.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.12%  "

# Row 33
$ws.Range("E33").Value = "  -4.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0741"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.74%  "

# Row 35
$ws.Range("E35").Value = "  -4.02%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.103"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.16%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.89%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0291"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.93%  "

# Row 40
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.05%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.71%  "

# Row 42
$ws.Range("E42").Value = "  -0.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.30%  "

# Row 44
$ws.Range("D44").Value = "1.806.24"
$ws.Range("E44").Value = "  +2.94%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.85%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.177"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.21%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.58%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.22%  "

# Row 49
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.95%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.97%  "

# Row 51
$ws.Range("D51").Value = "2.415.82"
$ws.Range("E51").Value = "  -2.17%  "
